# add date time parsing
# Rename header columns (id -> bio_id, data_type -> usage, data_source -> source,
# input_type -> type) across the "samples" and "references" sheets, normalise
# now-redundant cell styles back to the default style, tighten the vertical page
# margins, simplify the header/footer text, and fix up the active sheet/selection
# state left over from the previous edit session.

$wb = $excel.ActiveWorkbook

$wsSamples    = $wb.Worksheets.Item(1)   # "samples"
$wsReferences = $wb.Worksheets.Item(2)   # "references"
# (3rd sheet, "other", only needs the blanket page-setup changes applied below)

# --- samples sheet: rename header row -------------------------------------
$wsSamples.Range("A1").Value2 = "bio_id"
$wsSamples.Range("C1").Value2 = "usage"
$wsSamples.Range("D1").Value2 = "source"

# drop the redundant duplicate cell style (applyFont=true variant) from the
# "ncbi accession" cells so they fall back onto the default style index
$wsSamples.Range("C2:C4").NumberFormat = "General"

# --- references sheet: rename header row (B/C swap + A/D renames) ---------
$wsReferences.Range("A1").Value2 = "bio_id"
$wsReferences.Range("B1").Value2 = "usage"
$wsReferences.Range("C1").Value2 = "type"
$wsReferences.Range("D1").Value2 = "source"

# drop the redundant duplicate cell style here too
$wsReferences.Range("C1").NumberFormat = "General"
$wsReferences.Range("D1").NumberFormat = "General"
$wsReferences.Range("D2").NumberFormat = "General"

# fix the lingering selection left on this sheet
$wsReferences.Range("D1").Select() | Out-Null

# --- page margins: 1.05277777777778in -> 1.025in (top/bottom) on all sheets
foreach ($ws in $wb.Worksheets) {
    $ws.PageSetup.TopMargin = 73.8
    $ws.PageSetup.BottomMargin = 73.8
    $ws.PageSetup.CenterHeader = "&A"
    $ws.PageSetup.CenterFooter = "Page &P"
}

# --- active sheet / tab selection: "samples" becomes the active tab -------
$wsSamples.Activate()
